$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("URRIOLA ARISMENDIZ INGRID MARYURI", 31),
    @("ALAMA NIMA CLARITZA MABEL", 29),
    @("AGURTO ORDINOLA LISBET JAQUELIN", 29),
    @("MANUEL LEUNARDO PRADO BAILON", 29),
    @("ALBIRENA GARCIA ANGEELO ALONSO", 29),
    @("MARYURI OJEDA VALLE", 28),
    @("CORDOVA CARMEN ANGIE NATALLY", 27),
    @("ROMAN GALECIO MARITZA DEL PILAR", 27),
    @("JUAREZ CARMEN PIERRE ALEXANDER", 25),
    @("RUIDIAS FRIAS MELISSA VICTORIA", 25),
    @("VEGA ROBLEDO FERNANDO ERNESTO", 25),
    @("ATOCHE PALACIOS LUIS ANGEL", 15),
    @("CARREÑO PALACIOS KATHERINE DE LOS MILAGROS", 15),
    @("BERNAOLA CARMEN ZUMIKO YASHURY", 15),
    @("61097774", 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $cellA = $ws.Cells.Item($row, 1)
    $name = $data[$i][0]
    if ($name -match '^\d+$') {
        # Keep numeric-looking empadronador names (e.g. a stray DNI) as text.
        $cellA.NumberFormat = "@"
        $cellA.Value = $name
    } else {
        $cellA.Value = $name
    }
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
